# Auto-generated edit script: shifts the weekly Perejil price rows
# down by one observation (two data rows per date) and appends the
# newest weekly reading at the top of the block (rows 56-57), pushing
# the oldest reading (previously rows 138-139) to new rows 140-141.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column D (Fecha) for rows 56-139 ---
$ws.Range("D56").Value = 44679
$ws.Range("D57").Value = 44679
$ws.Range("D58").Value = 44330
$ws.Range("D59").Value = 44330
$ws.Range("D60").Value = 44476
$ws.Range("D61").Value = 44476
$ws.Range("D62").Value = 44386
$ws.Range("D63").Value = 44386
$ws.Range("D64").Value = 44579
$ws.Range("D65").Value = 44579
$ws.Range("D66").Value = 44488
$ws.Range("D67").Value = 44488
$ws.Range("D68").Value = 44237
$ws.Range("D69").Value = 44237
$ws.Range("D70").Value = 44292
$ws.Range("D71").Value = 44292
$ws.Range("D72").Value = 44616
$ws.Range("D73").Value = 44616
$ws.Range("D74").Value = 44490
$ws.Range("D75").Value = 44490
$ws.Range("D76").Value = 44609
$ws.Range("D77").Value = 44609
$ws.Range("D78").Value = 44453
$ws.Range("D79").Value = 44453
$ws.Range("D80").Value = 44320
$ws.Range("D81").Value = 44320
$ws.Range("D82").Value = 44475
$ws.Range("D83").Value = 44475
$ws.Range("D84").Value = 44327
$ws.Range("D85").Value = 44327
$ws.Range("D86").Value = 44350
$ws.Range("D87").Value = 44350
$ws.Range("D88").Value = 44574
$ws.Range("D89").Value = 44574
$ws.Range("D90").Value = 44523
$ws.Range("D91").Value = 44523
$ws.Range("D92").Value = 44400
$ws.Range("D93").Value = 44400
$ws.Range("D94").Value = 44252
$ws.Range("D95").Value = 44252
$ws.Range("D96").Value = 44299
$ws.Range("D97").Value = 44299
$ws.Range("D98").Value = 44460
$ws.Range("D99").Value = 44460
$ws.Range("D100").Value = 44334
$ws.Range("D101").Value = 44334
$ws.Range("D102").Value = 44565
$ws.Range("D103").Value = 44565
$ws.Range("D104").Value = 44405
$ws.Range("D105").Value = 44405
$ws.Range("D106").Value = 44358
$ws.Range("D107").Value = 44358
$ws.Range("D108").Value = 44383
$ws.Range("D109").Value = 44383
$ws.Range("D110").Value = 44582
$ws.Range("D111").Value = 44582
$ws.Range("D112").Value = 44525
$ws.Range("D113").Value = 44525
$ws.Range("D114").Value = 44435
$ws.Range("D115").Value = 44435
$ws.Range("D116").Value = 44442
$ws.Range("D117").Value = 44442
$ws.Range("D118").Value = 44194
$ws.Range("D119").Value = 44194
$ws.Range("D120").Value = 44398
$ws.Range("D121").Value = 44398
$ws.Range("D122").Value = 44512
$ws.Range("D123").Value = 44512
$ws.Range("D124").Value = 44222
$ws.Range("D125").Value = 44222
$ws.Range("D126").Value = 44285
$ws.Range("D127").Value = 44285
$ws.Range("D128").Value = 44203
$ws.Range("D129").Value = 44203
$ws.Range("D130").Value = 44274
$ws.Range("D131").Value = 44274
$ws.Range("D132").Value = 44607
$ws.Range("D133").Value = 44607
$ws.Range("D134").Value = 44425
$ws.Range("D135").Value = 44425
$ws.Range("D136").Value = 44250
$ws.Range("D137").Value = 44250
$ws.Range("D138").Value = 44306
$ws.Range("D139").Value = 44306

# --- Update column O (Origen) for rows with a changed region ---
$ws.Range("O62").Value = "Región de Ñuble"
$ws.Range("O63").Value = "Región de Ñuble"
$ws.Range("O64").Value = "Región Metropolitana"
$ws.Range("O65").Value = "Región Metropolitana"
$ws.Range("O100").Value = "Región de Ñuble"
$ws.Range("O101").Value = "Región de Ñuble"
$ws.Range("O102").Value = "Región Metropolitana"
$ws.Range("O103").Value = "Región Metropolitana"
$ws.Range("O134").Value = "Región de Ñuble"
$ws.Range("O135").Value = "Región de Ñuble"
$ws.Range("O136").Value = "Región de Arica y Parinacota"
$ws.Range("O137").Value = "Región de Arica y Parinacota"

# --- Update column J (Volumen) for rows with a changed volume ---
$ws.Range("J114").Value = 200
$ws.Range("J115").Value = 100
$ws.Range("J116").Value = 300
$ws.Range("J117").Value = 150

# --- Append new rows 140 and 141 (copies of what used to be rows 138/139) ---
$ws.Range("A140").Value = 11
$ws.Range("B140").Value = "Vega Monumental Concepción"
$ws.Range("C140").Value = "Bíobío"
$ws.Range("D140").Value = 44160
$ws.Range("E140").Value = 8
$ws.Range("F140").Value = 100112044
$ws.Range("G140").Value = "Perejil"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 200
$ws.Range("K140").Value = 600
$ws.Range("L140").Value = 700
$ws.Range("M140").Value = 650
$ws.Range("N140").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O140").Value = "Región de Ñuble"
$ws.Range("P140").Value = 650
$ws.Range("Q140").Value = 1
$ws.Range("R140").Value = "Hortaliza"
$ws.Range("D140").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A141").Value = 11
$ws.Range("B141").Value = "Vega Monumental Concepción"
$ws.Range("C141").Value = "Bíobío"
$ws.Range("D141").Value = 44160
$ws.Range("E141").Value = 8
$ws.Range("F141").Value = 100112044
$ws.Range("G141").Value = "Perejil"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Segunda"
$ws.Range("J141").Value = 100
$ws.Range("K141").Value = 500
$ws.Range("L141").Value = 500
$ws.Range("M141").Value = 500
$ws.Range("N141").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O141").Value = "Región de Ñuble"
$ws.Range("P141").Value = 500
$ws.Range("Q141").Value = 1
$ws.Range("R141").Value = "Hortaliza"
$ws.Range("D141").NumberFormat = "YYYY-MM-DD HH:MM:SS"
